# Update odds values in Sheet1 per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.05
$ws.Range("L2").Value = 5.5
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 2.05
$ws.Range("W2").Value = 8.5
$ws.Range("AC2").Value = 15
$ws.Range("AJ2").Value = 17
$ws.Range("AM2").Value = 41
$ws.Range("AO2").Value = 7.5
$ws.Range("AQ2").Value = 21
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 51
$ws.Range("BA2").Value = 101

# Row 3
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 5
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.5
$ws.Range("Y3").Value = 9.5
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AF3").Value = 81
$ws.Range("AH3").Value = 9
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 401

# Row 5
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 3.8
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 4.75
$ws.Range("N5").Value = 8.5
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 1.75
$ws.Range("W5").Value = 6
$ws.Range("AB5").Value = 34
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 19
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 451
$ws.Range("AH5").Value = 9
$ws.Range("AI5").Value = 19
$ws.Range("AJ5").Value = 13
$ws.Range("AN5").Value = 3.75
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 201
$ws.Range("AV5").Value = 67
$ws.Range("BB5").Value = 301

# Row 6
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 21

# Row 8
$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 3.6
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 1.91
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.48
$ws.Range("X8").Value = 9
$ws.Range("AI8").Value = 17
$ws.Range("AQ8").Value = 41
$ws.Range("AW8").Value = 5.5
$ws.Range("AZ8").Value = 81
$ws.Range("BA8").Value = 126

# Row 9
$ws.Range("G9").Value = 3.25
$ws.Range("I9").Value = 2.3
$ws.Range("L9").Value = 3
$ws.Range("X9").Value = 15
$ws.Range("AA9").Value = 26
$ws.Range("AB9").Value = 34
$ws.Range("AI9").Value = 11

# Row 10
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.7
